# "Practice tasks and final revisions"
# Renames each task-order sheet (new timestamped run id) and rewrites the
# stimulus-file / condition values that make up each sheet's task order.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16509961746594186"
$ws.Range("B2").Value = "go_stims-16509961746193836.csv"
$ws.Range("B3").Value = "GNG_stims-1650996174643414.csv"
$ws.Range("B4").Value = "go_stims-1650996174643414.csv"
$ws.Range("B5").Value = "GNG_stims-16509961746594186.csv"

# --- Sheet 2: NB -----------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16509961763554153"
$ws.Range("B2").Value = "TB-16509961763313882.csv"
$ws.Range("B3").Value = "OB-1650996175347419.csv"
$ws.Range("B4").Value = "ZB-match_8-16509961749313807.csv"
$ws.Range("B5").Value = "ZB-match_1-16509961750273814.csv"
$ws.Range("B6").Value = "OB-16509961757873776.csv"
$ws.Range("B7").Value = "ZB-match_8-1650996175227389.csv"
$ws.Range("B8").Value = "OB-1650996175883421.csv"
$ws.Range("B9").Value = "TB-16509961760514116.csv"
$ws.Range("B10").Value = "TB-16509961759874115.csv"

# --- Sheet 3: RS -------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16509961763554153"
$ws.Range("B2").Value = "eyes open"
$ws.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL --------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16509961764114149"
$ws.Range("B2").Value = "MM_stims-16509961763794134.csv"
$ws.Range("B3").Value = "ZM_stims-16509961763554153.csv"
$ws.Range("B4").Value = "MM_stims-1650996176395416.csv"
$ws.Range("B5").Value = "ZM_stims-16509961763794134.csv"
$ws.Range("B6").Value = "MM_stims-16509961764114149.csv"
$ws.Range("B7").Value = "ZM_stims-1650996176395416.csv"

# --- Sheet 5: vSAT ---------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16509961764754164"
$ws.Range("B2").Value = "SAT_stims-16509961764274187.csv"
$ws.Range("B3").Value = "vSAT_stims-16509961764434164.csv"
$ws.Range("B4").Value = "vSAT_stims-16509961764594145.csv"
$ws.Range("B5").Value = "SAT_stims-16509961764114149.csv"
